$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 7, shifting existing rows 7-28 down to 8-29.
$ws.Rows.Item(7).Insert()

# Populate the new row 7 with the new weekly data point.
$ws.Cells.Item(7, 1).Value = 9
$ws.Cells.Item(7, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(7, 3).Value = "Metropolitana"
$ws.Cells.Item(7, 4).Value = 44677
$ws.Cells.Item(7, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(7, 5).Value = 13
$ws.Cells.Item(7, 6).Value = 100112035
$ws.Cells.Item(7, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(7, 8).Value = "Sin especificar"
$ws.Cells.Item(7, 9).Value = "Primera"
$ws.Cells.Item(7, 10).Value = 34
$ws.Cells.Item(7, 11).Value = 25000
$ws.Cells.Item(7, 12).Value = 26000
$ws.Cells.Item(7, 13).Value = 25500
$ws.Cells.Item(7, 14).Value = '$/malla 15 kilos'
$ws.Cells.Item(7, 15).Value = "Hijuelas"
$ws.Cells.Item(7, 16).Value = 1700
$ws.Cells.Item(7, 17).Value = 15
$ws.Cells.Item(7, 18).Value = "Hortaliza"
